# Add a "second" row of data to Sheet1:
#   A3 = "second" (a brand new shared string)
#   B3 = "change" (same text/shared string already used in B2)
# This mirrors B2's "change" value while introducing one new unique
# shared string ("second"), growing the sheet from A1:B2 to A1:B3 and
# moving the active selection to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "second"
$ws.Range("B3").Value = "change"

$ws.Range("B3").Select()
